$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("M2").ClearContents()
$ws.Range("H19").Value = 2382.8333
$ws.Range("J19").Value = 2819.3
$ws.Range("L19").Value = 2819.3
$ws.Range("N19").Value = -3169.3
$ws.Range("H76").Value = 4000
$ws.Range("J76").Value = 4000
$ws.Range("L76").Value = 4000
$ws.Range("N76").Value = -4630
$ws.Range("H79").Value = 4000
$ws.Range("J79").Value = 4000
$ws.Range("L79").Value = 4000
$ws.Range("N79").Value = -6184
$ws.Range("H92").Value = 3908.0293
$ws.Range("I92").Value = 3647.75
$ws.Range("K92").Value = 3647.75
$ws.Range("M92").Value = -2399.75
$ws.Range("H106").Value = 100976.7
$ws.Range("I106").Value = 112140.89
$ws.Range("K106").Value = 112140.89
$ws.Range("M106").Value = -111509.89
$ws.Range("H112").Value = 4648183.5
$ws.Range("I112").Value = 1466.3334
$ws.Range("J112").Value = 5809862.5
$ws.Range("K112").Value = 4399.0002
$ws.Range("L112").Value = 17429587.5
$ws.Range("M112").Value = -3291.0002
$ws.Range("N112").Value = -17431803.5
$ws.Range("H113").Value = 18126.953
$ws.Range("I113").Value = 18064.562
$ws.Range("K113").Value = 18064.562
$ws.Range("M113").Value = -14810.562
$ws.Range("H118").Value = 7143025
$ws.Range("I118").Value = 7143025
$ws.Range("J118").Value = 0
$ws.Range("K118").Value = 21429075
$ws.Range("L118").Value = 0
$ws.Range("M118").Value = -21427418
$ws.Range("N118").ClearContents()
$ws.Range("H132").Value = 3129.5
$ws.Range("I132").Value = 1293.9259
$ws.Range("K132").Value = 3881.7777
$ws.Range("M132").Value = -1351.7777
$ws.Range("H138").Value = 3406.573
$ws.Range("I138").Value = 1413.6072
$ws.Range("J138").Value = 4227.206
$ws.Range("K138").Value = 4240.821599999999
$ws.Range("L138").Value = 12681.618
$ws.Range("M138").Value = 899.1784000000007
$ws.Range("N138").Value = -22961.618
$ws.Range("H141").Value = 6369.533
$ws.Range("I141").Value = 6403.3076
$ws.Range("K141").Value = 19209.9228
$ws.Range("M141").Value = -14029.9228

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 11849.149
$ws.Range("I32").Value = 11235.822
$ws.Range("K32").Value = 11235.822
$ws.Range("M32").Value = -10948.822
$ws.Range("H35").Value = 0
$ws.Range("I35").Value = 0
$ws.Range("K35").Value = 0
$ws.Range("M35").ClearContents()
$ws.Range("H61").Value = 7226.722
$ws.Range("I61").Value = 3294.7693
$ws.Range("J61").Value = 17449.8
$ws.Range("K61").Value = 3294.7693
$ws.Range("L61").Value = 17449.8
$ws.Range("M61").Value = -3082.7693
$ws.Range("N61").Value = -17873.8
$ws.Range("H136").Value = 7226.722
$ws.Range("I136").Value = 3294.7693
$ws.Range("J136").Value = 17449.8
$ws.Range("K136").Value = 9884.3079
$ws.Range("L136").Value = 52349.39999999999
$ws.Range("M136").Value = -7334.3079
$ws.Range("N136").Value = -57449.39999999999
$ws.Range("H140").Value = 83333.336
$ws.Range("J140").Value = 83333.336
$ws.Range("L140").Value = 83333.336
$ws.Range("N140").Value = -93693.336

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2947.1052
$ws.Range("J20").Value = 10466.333
$ws.Range("L20").Value = 10466.333
$ws.Range("N20").Value = -10960.333
$ws.Range("H33").Value = 13749.5
$ws.Range("I33").Value = 13749.5
$ws.Range("K33").Value = 13749.5
$ws.Range("M33").Value = -13413.5
$ws.Range("H94").Value = 863.53845
$ws.Range("I94").Value = 870.7
$ws.Range("K94").Value = 870.7
$ws.Range("M94").Value = -419.7
$ws.Range("H105").Value = 2073.7
$ws.Range("I105").Value = 1998.5555
$ws.Range("J105").Value = 2750
$ws.Range("K105").Value = 1998.5555
$ws.Range("L105").Value = 2750
$ws.Range("M105").Value = -251.5554999999999
$ws.Range("N105").Value = -6244
$ws.Range("H107").Value = 2293.6843
$ws.Range("I107").Value = 1347.5
$ws.Range("K107").Value = 1347.5
$ws.Range("M107").Value = 572.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 170439.03
$ws.Range("J31").Value = 40607.45
$ws.Range("L31").Value = 40607.45
$ws.Range("N31").Value = -41197.45
$ws.Range("H34").Value = 170439.03
$ws.Range("J34").Value = 40607.45
$ws.Range("L34").Value = 40607.45
$ws.Range("N34").Value = -41011.45
$ws.Range("H107").Value = 3830.1765
$ws.Range("J107").Value = 4772.769
$ws.Range("L107").Value = 4772.769
$ws.Range("N107").Value = -8612.769
$ws.Range("H132").Value = 3188.2341
$ws.Range("I132").Value = 2217.442
$ws.Range("K132").Value = 6652.326
$ws.Range("M132").Value = -4122.326
$ws.Range("H134").Value = 3429.291
$ws.Range("I134").Value = 3627.568
$ws.Range("K134").Value = 10882.704
$ws.Range("M134").Value = -8347.704000000002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H47").Value = 14222.923
$ws.Range("I47").Value = 9412.166999999999
$ws.Range("K47").Value = 28236.501
$ws.Range("M47").Value = -27805.501
$ws.Range("H55").Value = 8687.375
$ws.Range("J55").Value = 9450
$ws.Range("L55").Value = 28350
$ws.Range("N55").Value = -28704
$ws.Range("H99").Value = 5118.1113
$ws.Range("I99").Value = 4732.875
$ws.Range("J99").Value = 8200
$ws.Range("K99").Value = 14198.625
$ws.Range("L99").Value = 24600
$ws.Range("M99").Value = -11952.625
$ws.Range("N99").Value = -29092
$ws.Range("H109").Value = 1236.75
$ws.Range("I109").Value = 1236.75
$ws.Range("K109").Value = 3710.25
$ws.Range("M109").Value = -2670.25
$ws.Range("H120").Value = 3250
$ws.Range("I120").Value = 3250
$ws.Range("K120").Value = 9750
$ws.Range("M120").Value = -4912
$ws.Range("H122").Value = 399.16666
$ws.Range("J122").Value = 356.5
$ws.Range("L122").Value = 3208.5
$ws.Range("N122").Value = -8108.5
$ws.Range("H129").Value = 23811676
$ws.Range("I129").Value = 47620820
$ws.Range("J129").Value = 2532.4285
$ws.Range("K129").Value = 142862460
$ws.Range("L129").Value = 7597.2855
$ws.Range("M129").Value = -142857460
$ws.Range("N129").Value = -17597.2855
$ws.Range("H131").Value = 13890721
$ws.Range("I131").Value = 250000510
$ws.Range("J131").Value = 1909.0883
$ws.Range("K131").Value = 750001530
$ws.Range("L131").Value = 5727.2649
$ws.Range("M131").Value = -749996490
$ws.Range("N131").Value = -15807.2649
$ws.Range("H133").Value = 4075.7144
$ws.Range("J133").Value = 3500
$ws.Range("L133").Value = 10500
$ws.Range("N133").Value = -20620

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 15000
$ws.Range("I46").Value = 13333.333
$ws.Range("K46").Value = 13333.333
$ws.Range("M46").Value = -13177.333
$ws.Range("H93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").ClearContents()
$ws.Range("H126").Value = 5615.76
$ws.Range("I126").Value = 4275.3335
$ws.Range("K126").Value = 12826.0005
$ws.Range("M126").Value = -10356.0005
$ws.Range("H132").Value = 31030.316
$ws.Range("I132").Value = 37443.863
$ws.Range("J132").Value = 10364.444
$ws.Range("K132").Value = 112331.589
$ws.Range("L132").Value = 31093.332
$ws.Range("M132").Value = -109801.589
$ws.Range("N132").Value = -36153.33199999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").ClearContents()
$ws.Range("H136").Value = 4375.6665
$ws.Range("I136").Value = 4297.625
$ws.Range("J136").Value = 5000
$ws.Range("K136").Value = 12892.875
$ws.Range("L136").Value = 15000
$ws.Range("M136").Value = -10342.875
$ws.Range("N136").Value = -20100

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H63").Value = 24555.25
$ws.Range("J63").Value = 24998.334
$ws.Range("L63").Value = 24998.334
$ws.Range("N63").Value = -26246.334
$ws.Range("H66").Value = 24555.25
$ws.Range("J66").Value = 24998.334
$ws.Range("L66").Value = 74995.00199999999
$ws.Range("N66").Value = -81235.00199999999
$ws.Range("H81").Value = 3075.9333
$ws.Range("I81").Value = 3557.6667
$ws.Range("J81").Value = 1149
$ws.Range("K81").Value = 7115.3334
$ws.Range("L81").Value = 2298
$ws.Range("M81").Value = -6054.3334
$ws.Range("N81").Value = -4420
$ws.Range("H84").Value = 3075.9333
$ws.Range("I84").Value = 3557.6667
$ws.Range("J84").Value = 1149
$ws.Range("K84").Value = 35576.667
$ws.Range("L84").Value = 11490
$ws.Range("M84").Value = -30272.667
$ws.Range("N84").Value = -22098
$ws.Range("H96").Value = 65872
$ws.Range("I96").Value = 253998.5
$ws.Range("J96").Value = 3163.1667
$ws.Range("K96").Value = 253998.5
$ws.Range("L96").Value = 3163.1667
$ws.Range("M96").Value = -252625.5
$ws.Range("N96").Value = -5909.1667
$ws.Range("H107").Value = 546.4545000000001
$ws.Range("I107").Value = 427.14285
$ws.Range("K107").Value = 1281.42855
$ws.Range("M107").Value = 638.5714499999999
$ws.Range("H126").Value = 2062.75
$ws.Range("I126").Value = 1341.1818
$ws.Range("K126").Value = 4023.5454
$ws.Range("M126").Value = -1553.5454
$ws.Range("H132").Value = 1884.8677
$ws.Range("I132").Value = 883.381
$ws.Range("J132").Value = 2332.3403
$ws.Range("K132").Value = 2650.143
$ws.Range("L132").Value = 6997.0209
$ws.Range("M132").Value = -120.143
$ws.Range("N132").Value = -12057.0209
